# "Màj Journal de bord" - add the new log entry for Alexandre's sheet
# and move the active tab/selection from "Brian" to "Alexandre".

$wb = $excel.ActiveWorkbook

$wsAlexandre = $wb.Worksheets.Item("Alexandre")
$wsBrian     = $wb.Worksheets.Item("Brian")

# --- Append the new "Log" table row on the Alexandre sheet -----------------
$table = $wsAlexandre.ListObjects.Item(1)
$table.ListRows.Add() | Out-Null

# Copy the formatting of the previous data row (B7:C7) onto the new row
# (B8:C8) so the new cells pick up the same cell styles ("Date Column" /
# "Event Column") as every other row in the log, then fill in the values.
$wsAlexandre.Range("B7:C7").Copy()
$wsAlexandre.Range("B8:C8").PasteSpecial(-4122)   # xlPasteFormats
$wsAlexandre.Application.CutCopyMode = 0

$wsAlexandre.Cells.Item(8, 2).Value = 43158
$wsAlexandre.Cells.Item(8, 3).Value = "création d'une ébauche JSON"

# Keep the row height consistent with the rest of the table.
$wsAlexandre.Rows.Item(8).RowHeight = 30.75

# --- Move the active sheet / selection from Brian to Alexandre -------------
$wsAlexandre.Activate()
$wsAlexandre.Range("C10").Select()

$wsBrian.Range("C11").Select()
$wsAlexandre.Activate()
